# Auto-generated edit script applying the Behemoth_Profits market-data refresh.
# Updates static currentAveragePrice / LevePrice* / LeveProfit* columns (H:N)
# per worksheet, matching the scheduled-runner data refresh in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 74103.60000000001
$ws.Range("I21").Value = 66506.875
$ws.Range("K21").Value = 66506.875
$ws.Range("M21").Value = -66038.875
$ws.Range("H23").Value = 74103.60000000001
$ws.Range("I23").Value = 66506.875
$ws.Range("K23").Value = 66506.875
$ws.Range("M23").Value = -66272.875
$ws.Range("H105").Value = 87371.45
$ws.Range("J105").Value = 87371.45
$ws.Range("L105").Value = 87371.45
$ws.Range("N105").Value = -94359.45
$ws.Range("H137").Value = 7665.1924
$ws.Range("I137").Value = 1546.125
$ws.Range("K137").Value = 4638.375
$ws.Range("M137").Value = -2088.375
$ws.Range("H138").Value = 3777.3118
$ws.Range("J138").Value = 4034.4443
$ws.Range("L138").Value = 12103.3329
$ws.Range("N138").Value = -22383.3329
$ws.Range("H141").Value = 1786.9166
$ws.Range("I141").Value = 1272.5555
$ws.Range("K141").Value = 3817.6665
$ws.Range("M141").Value = 1362.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13910948
$ws.Range("I32").Value = 17871760
$ws.Range("K32").Value = 17871760
$ws.Range("M32").Value = -17871473
$ws.Range("H61").Value = 62635132
$ws.Range("I61").Value = 125005670
$ws.Range("J61").Value = 264589.25
$ws.Range("K61").Value = 125005670
$ws.Range("L61").Value = 264589.25
$ws.Range("M61").Value = -125005458
$ws.Range("N61").Value = -265013.25
$ws.Range("H88").Value = 2840.5
$ws.Range("J88").Value = 1825.75
$ws.Range("L88").Value = 1825.75
$ws.Range("N88").Value = -2637.75
$ws.Range("H91").Value = 2840.5
$ws.Range("J91").Value = 1825.75
$ws.Range("L91").Value = 1825.75
$ws.Range("N91").Value = -4633.75
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 62635132
$ws.Range("I136").Value = 125005670
$ws.Range("J136").Value = 264589.25
$ws.Range("K136").Value = 375017010
$ws.Range("L136").Value = 793767.75
$ws.Range("M136").Value = -375014460
$ws.Range("N136").Value = -798867.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 41067.75
$ws.Range("J63").Value = 21423.666
$ws.Range("L63").Value = 21423.666
$ws.Range("N63").Value = -22795.666
$ws.Range("H66").Value = 41067.75
$ws.Range("J66").Value = 21423.666
$ws.Range("L66").Value = 64270.99800000001
$ws.Range("N66").Value = -71134.99800000001
$ws.Range("H94").Value = 864.4583
$ws.Range("J94").Value = 1701.6666
$ws.Range("L94").Value = 1701.6666
$ws.Range("N94").Value = -2603.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 981501.3
$ws.Range("I31").Value = 18858.334
$ws.Range("K31").Value = 18858.334
$ws.Range("M31").Value = -18563.334
$ws.Range("H34").Value = 981501.3
$ws.Range("I34").Value = 18858.334
$ws.Range("K34").Value = 18858.334
$ws.Range("M34").Value = -18656.334
$ws.Range("H55").Value = 1573
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H81").Value = 122927.5
$ws.Range("J81").Value = 122927.5
$ws.Range("L81").Value = 122927.5
$ws.Range("N81").Value = -124923.5
$ws.Range("H84").Value = 122927.5
$ws.Range("J84").Value = 122927.5
$ws.Range("L84").Value = 368782.5
$ws.Range("N84").Value = -378766.5
$ws.Range("H99").Value = 1908.8
$ws.Range("I99").Value = 1811.0667
$ws.Range("J99").Value = 2202
$ws.Range("K99").Value = 1811.0667
$ws.Range("L99").Value = 2202
$ws.Range("M99").Value = -313.0667000000001
$ws.Range("N99").Value = -5198
$ws.Range("H126").Value = 1908.8
$ws.Range("I126").Value = 1811.0667
$ws.Range("J126").Value = 2202
$ws.Range("K126").Value = 5433.2001
$ws.Range("L126").Value = 6606
$ws.Range("M126").Value = -2963.2001
$ws.Range("N126").Value = -11546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3900
$ws.Range("I58").Value = 2800
$ws.Range("K58").Value = 8400
$ws.Range("M58").Value = -8272
$ws.Range("H131").Value = 6784.727
$ws.Range("I131").Value = 8441.875
$ws.Range("J131").Value = 2365.6667
$ws.Range("K131").Value = 25325.625
$ws.Range("L131").Value = 7097.000100000001
$ws.Range("M131").Value = -20285.625
$ws.Range("N131").Value = -17177.0001
$ws.Range("H132").Value = 2631.3635
$ws.Range("I132").Value = 2421.2856
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 21791.5704
$ws.Range("L132").Value = 26991
$ws.Range("M132").Value = -19261.5704
$ws.Range("N132").Value = -32051

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5780
$ws.Range("I70").Value = 5666
$ws.Range("K70").Value = 5666
$ws.Range("M70").Value = -5396
$ws.Range("H73").Value = 5780
$ws.Range("I73").Value = 5666
$ws.Range("K73").Value = 5666
$ws.Range("M73").Value = -4730
$ws.Range("H102").Value = 2958.6667
$ws.Range("I102").Value = 2454.7083
$ws.Range("K102").Value = 2454.7083
$ws.Range("M102").Value = -832.7082999999998
$ws.Range("H122").Value = 3226.125
$ws.Range("I122").Value = 3072.4736
$ws.Range("J122").Value = 3810
$ws.Range("K122").Value = 9217.4208
$ws.Range("L122").Value = 11430
$ws.Range("M122").Value = -6767.4208
$ws.Range("N122").Value = -16330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5403.8335
$ws.Range("I22").Value = 5403.8335
$ws.Range("K22").Value = 5403.8335
$ws.Range("M22").Value = -5108.8335
$ws.Range("H27").Value = 5403.8335
$ws.Range("I27").Value = 5403.8335
$ws.Range("K27").Value = 5403.8335
$ws.Range("M27").Value = -5296.8335
$ws.Range("H40").Value = 4230.76
$ws.Range("I40").Value = 3938.2778
$ws.Range("K40").Value = 3938.2778
$ws.Range("M40").Value = -3802.2778
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H55").Value = 83333930
$ws.Range("J55").Value = 706.8
$ws.Range("L55").Value = 706.8
$ws.Range("N55").Value = -1052.8
$ws.Range("H106").Value = 24221.125
$ws.Range("J106").Value = 24221.125
$ws.Range("L106").Value = 24221.125
$ws.Range("N106").Value = -26745.125
$ws.Range("H122").Value = 5942.154
$ws.Range("I122").Value = 5356.952
$ws.Range("K122").Value = 16070.856
$ws.Range("M122").Value = -13620.856
$ws.Range("H125").Value = 129000
$ws.Range("J125").Value = 129000
$ws.Range("L125").Value = 129000
$ws.Range("N125").Value = -138840
$ws.Range("H136").Value = 355418.84
$ws.Range("I136").Value = 337502.66
$ws.Range("J136").Value = 373335
$ws.Range("K136").Value = 1012507.98
$ws.Range("L136").Value = 1120005
$ws.Range("M136").Value = -1009957.98
$ws.Range("N136").Value = -1125105

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 38888
$ws.Range("I18").Value = 38888
$ws.Range("K18").Value = 38888
$ws.Range("M18").Value = -38715
$ws.Range("H29").Value = 37249.25
$ws.Range("I29").Value = 31999
$ws.Range("J29").Value = 42499.5
$ws.Range("K29").Value = 31999
$ws.Range("L29").Value = 42499.5
$ws.Range("M29").Value = -31709
$ws.Range("N29").Value = -43079.5
$ws.Range("H61").Value = 24918.223
$ws.Range("I61").Value = 12852
$ws.Range("J61").Value = 40001
$ws.Range("K61").Value = 12852
$ws.Range("L61").Value = 40001
$ws.Range("M61").Value = -12560
$ws.Range("N61").Value = -40585
$ws.Range("H63").Value = 26075.334
$ws.Range("H64").Value = 49993.5
$ws.Range("J64").Value = 49993.5
$ws.Range("L64").Value = 49993.5
$ws.Range("N64").Value = -50489.5
$ws.Range("H66").Value = 26075.334
$ws.Range("H67").Value = 49993.5
$ws.Range("J67").Value = 49993.5
$ws.Range("L67").Value = 49993.5
$ws.Range("N67").Value = -51709.5
$ws.Range("H104").Value = 72775.8
$ws.Range("J104").Value = 72775.8
$ws.Range("L104").Value = 72775.8
$ws.Range("N104").Value = -79763.8
$ws.Range("H122").Value = 4778.4736
$ws.Range("I122").Value = 3399.5293
$ws.Range("J122").Value = 16499.5
$ws.Range("K122").Value = 10198.5879
$ws.Range("L122").Value = 49498.5
$ws.Range("M122").Value = -7748.5879
$ws.Range("N122").Value = -54398.5
